# Update column G ("K") values for rows 2-9 on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2 = 1
    3 = 1
    4 = 4
    5 = 1
    6 = 1
    7 = 4
    8 = 2
    9 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
